$d = $word.ActiveDocument

# The paragraph currently reads (with a hidden "_GoBack" bookmark placed
# right before "user base"):
#   ... potential for the rapid growth of the [bookmark]user base which can
#   be further used advertisement revenue models.
#
# It needs to become:
#   ... potential for the rapid growth of the user base which can
#   be further used for [bookmark]advertisement revenue models.
#
# i.e. the word "for " is inserted right before "advertisement", and the
# "_GoBack" bookmark moves from before "user base" to right before
# "advertisement". Locate both anchor points via Find so the edit is not
# dependent on hard-coded character offsets.

# --- locate "advertisement revenue models." ---
$findAd = $d.Content
$findAd.Find.Execute("advertisement revenue models.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$adStart = $findAd.Start

# --- locate the run boundary right after the lone "b" run, i.e. right
#     after "user base which can b" ---
$findB = $d.Content
$findB.Find.Execute("user base which can b", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterB = $findB.End

# Temporary "barrier" bookmarks keep the surrounding runs from being
# coalesced by the engine's run-normalisation pass when we edit the text
# in between them (mirrors the paragraph's existing pattern of runs being
# split around the "_GoBack" bookmark).
$d.Bookmarks.Add("zzTmpBarrier1", $d.Range($afterB, $afterB))

# Insert the missing word right before "advertisement".
$insertedWord = "for "
$insertPoint = $d.Range($adStart, $adStart)
$insertPoint.InsertBefore($insertedWord)

# After the insertion, "advertisement" (and everything from $adStart on)
# has shifted right by the length of the inserted text.
$newAdStart = $adStart + $insertedWord.Length

# Barrier that separates "e further used " from "for ".
$d.Bookmarks.Add("zzTmpBarrier2", $d.Range($adStart, $adStart))

# Move "_GoBack" so it sits right before "advertisement" again.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($newAdStart, $newAdStart))

# Clean up the scaffolding bookmarks - only "_GoBack" should remain.
$d.Bookmarks.Item("zzTmpBarrier1").Delete()
$d.Bookmarks.Item("zzTmpBarrier2").Delete()
